$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing extr1..extr8 rows (currently rows 8..15) down two rows, to
# rows 10..17, to make room for the two new "line7"/"line8" rows. Copy cell by
# cell (bottom-up) instead of using Rows.Insert() so we don't create a spare
# duplicate style record for the inherited row formatting.
for ($r = 15; $r -ge 8; $r--) {
    $dest = $r + 2
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($dest, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# --- New rows 8 and 9: line7 / line8 ---
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Apply the same formatting (bold, centered, bordered) used by the other
# column-A index cells to the two new rows, and to rows 16/17 which are past
# the original used range and so start out with no style at all.
$styleRows = @(8, 9, 16, 17)
foreach ($r in $styleRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- Re-number / refill the rows that used to be extr1..extr8 (now rows 10..17) ---
# row 10 = extr1
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# row 11 = extr2
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# row 12 = extr3
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $false

# row 13 = extr4
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $true

# row 14 = extr5
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

# row 15 = extr6
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

# row 16 = extr7
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

# row 17 = extr8
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false
